$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so values like "245.04" are not
# auto-converted to numbers by Excel's type inference, matching the
# original inline-string storage. Style is restored to Normal afterwards
# so no visible formatting changes are introduced.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "36.346.89"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").Value = "1.981.20"
$ws.Range("E3").Value = "  -3.96%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "245.04"
$ws.Range("E5").Value = "  -3.47%  "
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("D7").Value = "58.99"
$ws.Range("E7").Value = "  -12.81%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.374"
$ws.Range("E9").Value = "  -4.86%  "
$ws.Range("D10").Value = "57.15"
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("D11").Value = "0.0831"
$ws.Range("E11").Value = "  +7.49%  "
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "23.26"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "0.860"
$ws.Range("E14").Value = "  -8.42%  "
$ws.Range("D15").Value = "13.96"
$ws.Range("E15").Value = "  -7.40%  "
$ws.Range("D16").Value = "2.271.65"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").Value = "5.45"
$ws.Range("E17").Value = "  -4.06%  "
$ws.Range("D18").Value = "1.979.62"
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("D19").Value = "36.184.61"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").Value = "70.30"
$ws.Range("E20").Value = "  -4.79%  "
$ws.Range("D21").Value = "0.0₃0875"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "5.29"
$ws.Range("E22").Value = "  -4.07%  "
$ws.Range("D23").Value = "233.86"
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  -5.89%  "
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  -6.49%  "
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("D28").Value = "162.45"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "0.132"
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("D30").Value = "19.78"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").Value = "4.87"
$ws.Range("E33").Value = "  -7.40%  "
$ws.Range("D34").Value = "0.0681"
$ws.Range("E34").Value = "  +7.39%  "
$ws.Range("D35").Value = "4.39"
$ws.Range("E35").Value = "  -7.87%  "
$ws.Range("D36").Value = "6.20"
$ws.Range("E36").Value = "  -2.23%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "2.25"
$ws.Range("D39").Value = "1.81"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  -5.35%  "
$ws.Range("E41").Value = "  -4.17%  "
$ws.Range("D42").Value = "0.0962"
$ws.Range("E42").Value = "  -7.08%  "
$ws.Range("E43").Value = "  -5.44%  "
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("E45").Value = "  -5.91%  "
$ws.Range("D46").Value = "16.15"
$ws.Range("E46").Value = "  -12.32%  "
$ws.Range("D47").Value = "92.05"
$ws.Range("E47").Value = "  -5.81%  "
$ws.Range("D50").Value = "2.82"
$ws.Range("E50").Value = "  -4.47%  "
$ws.Range("D51").Value = "45.05"
$ws.Range("E51").Value = "  -4.59%  "

# Row 48/49: FraxShare and Maker swap rank positions with updated price/volume.
$ws.Range("D48").Value = "1.363.78"
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("D49").Value = "7.44"
$ws.Range("E49").Value = "  -7.06%  "

# Restore default (Normal) style on column D now that the text values are set.
$colD.Style = "Normal"

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
